$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# The cell "适用于元素快速插入和删除" + bookmark(_GoBack) + "的场景。" becomes a
# single run "适用于元素快速插入和删除的场景。" with the bookmark removed.
# A Find/Replace whose search span crosses the run/bookmark boundary makes
# Word coalesce the text into one run and drop the now-redundant bookmark.
$null = $d.Content.Find.Execute(
    "适用于元素快速插入和删除的场景。", $true, $false, $false, $false, $false,
    $true, 1, $false, "适用于元素快速插入和删除的场景。", 2)

# --- Change 2 -------------------------------------------------------------
# The cell "需要排序时可用。需要快速查找可用。不能制定元素插入位置。" (single run)
# becomes two runs: "...不能指定" + bookmark(_GoBack) + "元素插入位置。"
# (note the "制定" -> "指定" typo fix baked into the text change).

# 2a) fix the typo first, staying inside the same run.
$null = $d.Content.Find.Execute(
    "不能制定", $true, $false, $false, $false, $false,
    $true, 1, $false, "不能指定", 2)

# 2b) locate the split point (just before "元素插入位置。") and drop a
# bookmark there, splitting the run in two exactly like the diff shows.
$r = $d.Content
$null = $r.Find.Execute(
    "元素插入位置。", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
$r.Collapse(1)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
